# Applies the "Penalty Reward System" forecast refresh edit:
#  - Shifts each weekly row's Week_Start_Date forward by one week
#  - Updates the MyForecast (column D) values on the "Forecast Comparison" sheet
#  - Refreshes the computed statistics on the "Summary" sheet

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: Week_Start_Date (B) & MyForecast (D) ---
# Each row's date is pushed one week later, and MyForecast gets a new value.
$forecastRows = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 72 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 72 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 73 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 74 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 74 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 74 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 73 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 74 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 75 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 76 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 76 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 76 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 76 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 69 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 69 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 69 }
)

foreach ($item in $forecastRows) {
    $r = $item.Row
    # Prefix with an apostrophe so the date-like text is kept as plain text
    # instead of being auto-converted into a date serial number.
    $wsForecast.Cells.Item($r, 2).Value = "'" + $item.Date
    $wsForecast.Cells.Item($r, 4).Value = $item.Forecast
}

# --- Summary sheet: refreshed statistics ---
$wsSummary.Range("B2").Value  = "2024-02-04 to 2025-01-05"
$wsSummary.Range("B4").Value  = "'128"
$wsSummary.Range("B5").Value  = "'46"
$wsSummary.Range("B6").Value  = "'38"
$wsSummary.Range("B8").Value  = "2193 units"
$wsSummary.Range("B9").Value  = "'1172"
$wsSummary.Range("B10").Value = "'585"
$wsSummary.Range("B11").Value = "'290"
$wsSummary.Range("B12").Value = "'76"
$wsSummary.Range("B13").Value = "'2025-04-06"
$wsSummary.Range("B14").Value = "'69"
$wsSummary.Range("B15").Value = "'2025-04-27"
